$d = $word.ActiveDocument

# 1. Date change in the first line.
$d.Content.Find.Execute("06.08.24", $false, $false, $false, $false, $false, $true, 1, $false, "05.08.24", 2)

# 2. Title paragraph: replace the old title text and drop the trailing line break
#    by overwriting the whole title paragraph's range (mark excluded).
$newTitle = "Improving Text Embeddings for Smaller Language Models Using Contrastive Fine-tuning"
$d.Paragraphs(2).Range.Text = $newTitle

# 3. Intro paragraph.
$old3 = "חוזרים לסקור מאמרים על מודלי דיפוזיה עם מאמר כחול לבן של קבוצת חוקרים מאוניברסיטת תל אביב. הם מציעים שיטה מעניינת לעריכה מהירה של תמונה. כלומר בהינתן תמונה עם פרומפט נתון c אנו רוצים ליצור תמונה עם פרומפט אחר c1."
$new3 = "חוזרים לסקור מאמרים קלילים על מודלי שפה והיום בפוקוס מודלי שפה קטנים. המאמר שנסקור קצרות היום מציע שיטה לשיפור ייצוג של טקסט המופק על ידי מודל שפה קטן. ידוע שמודל שפה קטן (במאמר שיפרו את הייצוגים של הדקודרים) לא תמיד מצטיין ביצירה של ייצוג (אמבדינג) עוצמתי של טקסט - פשוט בגלל הגודל ו-expressiveness נמוכה יחסית."
$d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# 4. Contrastive learning explanation paragraph.
$old4 = "כמו שאתם זוכרים מודלי דיפוזיה מגנרטים תמונה על ידי הסרה רעש הדרגתית (denoising). בכל שלב המודל חוזה כמה רעש צריך להסיר מהתמונה והרעש המשוערך הזה מחוסר מהתמונה המורעשת באיטרציה הקודמת. השיטה הפשוטה לעשות עריכה של תמונה היא:"
$new4 = "אז המאמר מציע להשתמש בשיטת למידה ניגודית (contrastive learning) כדי לשפר את הביצועים. בגדול למידה ניגודית מאמנת מודל (לייצוג דאטה) במטרה לקרב פיסות דאטה (למשל תמונות או טקסט) שהן קרובות (סמנטית או בעלות אותה משמעות) ובאותו הזמן להרחיק את הייצוגים של פיסות דאטה לא דומות. השיטה הוצגה ב- 2018 על ידי Oord האגדי. "
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2)

# 5. Fine-tuning with LoRA paragraph.
$old5 = "להחסיר מהתמונה(המקורית) באיטרציה t את הרעש הזה המשוערך עם פרומפט c (כמו שעושים כאשר אין עריכה) "
$new5 = "המאמר מציע להשתמש בלמידה ניגודית כדי לעשות פיין טיון לייצוגי הדאטה המופקים על ידי מודל שפה בפרט הפלט של השכבה האחרונה עבור טוקן EoS המסמן את סוף המשפט. עדכון משקלי המודל נעשה כמובן עם LoRA על דאטהסט המכיל משפטים בעלי משמעות קרובה וגם זוגות משפטים רחוקים סמנטית. המחברים טוענים שזה משפר את איכות הייצוג המופק על ידי המודל למספר משימות downstream (בפרט סיווג)."
$d.Content.Find.Execute($old5, $false, $false, $false, $false, $false, $true, 1, $false, $new5, 2)

# 6. Closing remark paragraph.
$old6 = "להוסיף אל התוצאה את התוחלת המשוערכת של התמונה המורעשת(הערוכה) עם הפרומפט c1 החדש (עם התמונה המורעשת הערוכה. "
$new6 = "מאמר קלילי ונעים לקריאה…."
$d.Content.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $new6, 2)

# 7. Remove the now-obsolete paragraphs (old paragraphs 8-13), keeping paragraph 7's mark,
#    then turn what remains of paragraph 7 into the new link line.
$delStart = $d.Paragraphs(8).Range.Start
$delEnd = $d.Paragraphs(13).Range.End
$d.Range($delStart, $delEnd).Delete()

$d.Paragraphs(7).Range.Text = "https://arxiv.org/abs/2408.00690 "
